$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# Update title text from "Video?" to "Video"
$titleShape = $s.Shapes.Item("Title 1")
$titleShape.TextFrame.TextRange.Text = "Video"

# Add the YouTube link (as a hyperlinked run, followed by a trailing space
# run) to the previously-empty content placeholder below the title.
$url = "https://www.youtube.com/watch?v=oHNqW97JLAo"
$contentShape = $s.Shapes.Item("Content Placeholder 2")
$tr = $contentShape.TextFrame.TextRange
$tr.Text = $url + " "
$linkRange = $tr.Characters(1, $url.Length)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = $url
